# Fix doubled email in signature blocks.
#
# The "Email:" paragraph correctly shows {{attorney_email}} inside a
# <w:hyperlink> run. A stray, duplicate standalone paragraph containing
# only a tab + " {{attorney_email}}" (no hyperlink) was also left in the
# document, causing the attorney's email address to render twice in the
# signature block. Remove that duplicate paragraph entirely, leaving the
# hyperlinked copy untouched.

$d = $word.ActiveDocument

# Walk paragraphs back-to-front (so deleting doesn't shift the indices of
# paragraphs we still need to examine) and drop any paragraph whose text
# still contains the {{attorney_email}} placeholder but which does NOT
# itself carry a hyperlink -- i.e. the duplicate plain-text copy, not the
# legitimate "Email:" line.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $rng = $p.Range
    if ($rng.Text -like "*{{attorney_email}}*" -and $rng.Hyperlinks.Count -eq 0) {
        $rng.Delete()
    }
}
